$d = $word.ActiveDocument

# 1. Update the title heading in the blue banner table
$d.Content.Find.Execute("Circle Language Spec: Interfaces", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Circle Language Construct Drafts | Interfaces", 2)

# 2. "In Text Code" -> "i" + "n Text Code" (lower-case the leading I, split into two runs)
$rng = $d.Content
$found = $rng.Find.Execute("In Text Code", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)
if ($found) {
    $start = $rng.Start

    # Change the leading "I" to lower-case "i"
    $iChar = $d.Range($start, $start + 1)
    $iChar.Text = "i"

    # Split the run between "i" and "n Text Code" by inserting then
    # removing a paragraph mark at that position.
    $splitPoint1 = $d.Range($start + 1, $start + 1)
    $splitPoint1.InsertParagraphAfter()
    $markPos1 = $d.Range($start + 1, $start + 2)
    $markPos1.Delete()

    # Split the run between "Interface Assignment " and "i" the same way.
    $splitPoint2 = $d.Range($start, $start)
    $splitPoint2.InsertParagraphAfter()
    $markPos2 = $d.Range($start, $start + 1)
    $markPos2.Delete()
}

# 3. Add a trailing period after "It can be pronounced 'interface is pointer to'"
$quoteOpen = [char]0x2018
$quoteClose = [char]0x2019
$d.Content.Find.Execute("It can be pronounced " + $quoteOpen + "interface is pointer to" + $quoteClose,
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "It can be pronounced " + $quoteOpen + "interface is pointer to" + $quoteClose + ".", 2)

# 4. Remove the hidden _GoBack bookmark
try {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
} catch {
}
